$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
    }
}

# 1. "Texts, diagrams or pictures would be all fine" - collapse split runs / remove proofErr marks (no visible text change)
Replace-Text "Texts, diagrams or pictures would be all fine" "Texts, diagrams or pictures would be all fine"

# 2. "...implement anything...footage is a hand..." - collapse split runs (no visible text change)
Replace-Text " we didn't implement anything that will tell the program that the object in the footage is a hand. In fact, any object can be used, and the program will identify this as a hand. A work around was implemented. A square box is used and placed onto the screen. Anything inside of this box will be seen in the program as a hand. It is essentially up to the user to display a hand in this box." " we didn't implement anything that will tell the program that the object in the footage is a hand. In fact, any object can be used, and the program will identify this as a hand. A work around was implemented. A square box is used and placed onto the screen. Anything inside of this box will be seen in the program as a hand. It is essentially up to the user to display a hand in this box."

# 3. Remove stray "of" - "During testing of we found" -> "During testing we found"
Replace-Text "During testing of we found" "During testing we found"

# 4. "the original values can be set" -> "they can be set"
Replace-Text "However, the original values can be set to different values" "However, they can be set to different values"

# 5. "This is leads to the assumption...is a hand." -> "However, this leads to the assumption...is the hand."
Replace-Text "This is leads to the assumption that the biggest white blob in the image is a hand." "However, this leads to the assumption that the biggest white blob in the image is the hand."

# 6. "(i.e., one, two, three, four hand signs)" -> "(i.e., one, two, three, or four hand signs)"
Replace-Text "(i.e., one, two, three, four hand signs)" "(i.e., one, two, three, or four hand signs)"

# 7. ". So we counted" -> ". So, we counted" (also collapses split runs/proofErr)
Replace-Text ". So we counted" ". So, we counted"

# 8. "counted all the defects that were greater than this value." -> "...greater than four digits long."
Replace-Text "counted all the defects that were greater than this value." "counted all the defects that were greater than four digits long."

# 9. Reorder the white marker / green line sentence
Replace-Text ". So, a white marker shows the user where the defects are when it is detected, and a green line in between fingers." ". So, a white marker and a green line shows the user where the defects are when it is detected between fingers."

# 10. "Link to a 2-minute Youtube video" - collapse split runs / remove proofErr marks (no visible text change)
Replace-Text "Link to a 2-minute Youtube video" "Link to a 2-minute Youtube video"

# 11. "To use the yellow box must be placed" -> "To use, the yellow box must be placed"
Replace-Text "To use the yellow box must be placed" "To use, the yellow box must be placed"

# 12. "must be good otherwise the user will notice noise" -> "must be decent otherwise the user may notice noise"
Replace-Text "must be good otherwise the user will notice noise in the binary image." "must be decent otherwise the user may notice noise in the binary image."

# 13. "were good at the time of the program's inception" -> "were performed well at the time of the program's inception"
Replace-Text "The Trackbar has values set that were good at the time of the program’s inception." "The Trackbar has values set that were performed well at the time of the program’s inception."

# 14. "a range with Lower and Upper boundaries" -> "a range with lower and upper boundaries"
Replace-Text "a range with Lower and Upper boundaries." "a range with lower and upper boundaries."

# 15. "To end the program, simple press the key" -> "To end the program, simply press the key"
Replace-Text "To end the program, simple press the key" "To end the program, simply press the key"
